$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the TestAcc header to reflect the new 134-sign test, and add a new
# "Comments" column after it.
$ws.Range("O1").Value = "TestAcc (134 signs)"
$ws.Range("P1").Value = "Comments"

# New TestAcc (134 signs) results for the Stride-1 Chameleon runs.
$ws.Range("O5").Value = 0.2388
$ws.Range("O8").Value = 0.37
$ws.Range("O11").Value = 0.41
$ws.Range("O14").Value = 0.47

# Note for the last stride-1 run.
$ws.Range("P14").Value = "Analysis by sign and by subject"

# Highlight the BestValAcc cells tied to the new TestAcc measurements.
$ws.Range("N5").Interior.Color = 65535
$ws.Range("N8").Interior.Color = 65535
$ws.Range("N11").Interior.Color = 65535
$ws.Range("N14").Interior.Color = 65535

$ws.Range("O6").Select() | Out-Null
